# Bump the "Förändrad" (Changed) date in column C from 2023-09-02 (45171)
# to 2023-09-03 (45172) for every data row (rows 2-387).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C2:C387").Value2 = 45172
